# ------------------------------------------------------------------
# Supplementary_Table_04.xlsx edit:
#   1) Rename the sheet "Comparison with other CoVs" ->
#      "Count of ZAP-binding motifs"
#   2) Add 12 new columns (Z:AK) with ZAP-motif count data:
#      cs_01_zap_motifs_count ... cs_11_zap_motifs_count, total_zap_count
# ------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Rename the worksheet
$ws.Name = "Count of ZAP-binding motifs"

# 2) Give the new header cells (Z1:AK1) the same bold/bordered/centered
#    style already used by the existing header row, by copying formats
#    from an equal-width block of existing header cells (N1:Y1, 12 cols).
$ws.Range("N1:Y1").Copy()
$ws.Range("Z1:AK1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$newCols = @("Z","AA","AB","AC","AD","AE","AF","AG","AH","AI","AJ","AK")

$headers = @(
    "cs_01_zap_motifs_count",
    "cs_02_zap_motifs_count",
    "cs_03_zap_motifs_count",
    "cs_04_zap_motifs_count",
    "cs_05_zap_motifs_count",
    "cs_06_zap_motifs_count",
    "cs_07_zap_motifs_count",
    "cs_08_zap_motifs_count",
    "cs_09_zap_motifs_count",
    "cs_10_zap_motifs_count",
    "cs_11_zap_motifs_count",
    "total_zap_count"
)

for ($c = 0; $c -lt $newCols.Length; $c++) {
    $ws.Range($newCols[$c] + "1").Value = $headers[$c]
}

# Data rows 2-15 (sheet rows) for columns Z:AK, one array per row in the
# same left-to-right order as $newCols above. Last column is the
# row total (sum of the 11 conserved-stretch counts).
$data = @(
    @(1,0,0,0,0,0,0,0,0,0,0,1),
    @(0,0,0,0,1,0,0,0,0,0,0,1),
    @(0,0,0,0,0,0,0,0,0,0,0,0),
    @(4,1,0,0,0,0,0,0,0,0,0,5),
    @(4,1,0,0,0,0,0,0,0,0,0,5),
    @(2,0,1,0,0,0,0,1,0,0,0,4),
    @(0,2,4,1,0,0,0,1,0,1,0,9),
    @(2,1,5,1,0,2,0,1,0,0,0,12),
    @(3,0,1,0,0,0,0,1,0,0,0,5),
    @(0,0,2,0,0,0,0,1,2,0,0,5),
    @(1,0,0,0,0,0,0,0,0,0,0,1),
    @(2,0,1,0,0,0,0,0,1,0,0,4),
    @(0,1,2,1,0,0,0,3,0,0,0,7),
    @(0,2,0,4,0,2,0,0,0,1,1,10)
)

$firstDataRow = 2
for ($r = 0; $r -lt $data.Length; $r++) {
    $rowNum = $firstDataRow + $r
    $rowVals = $data[$r]
    for ($c = 0; $c -lt $newCols.Length; $c++) {
        $ws.Range($newCols[$c] + $rowNum).Value = $rowVals[$c]
    }
}
